# household_member.xlsx - "hide notes in contents screen" edit
#
# Adds a new `hideInContents` column to the `survey` sheet and marks the
# existing `note` row (row 6, type = "note") as TRUE so that note-type
# fields are hidden from the contents screen. Also restores focus/selection
# to the cell the user would naturally have ended on after typing the new
# column (H7, right below the new header/value pair), which becomes the
# active sheet/tab in the saved workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

# New column header, row 1 (H1): "hideInContents"
$ws.Range("H1").Value = "hideInContents"

# Flag the "note" row (row 6) to be hidden from the contents screen.
$ws.Range("H6").Value = $true

# Give the new column a sensible width, same as a user dragging/auto-sizing
# the column border after adding the new field.
$ws.Columns.Item(8).ColumnWidth = 14

# Make "survey" the active sheet/tab and leave the selection on the next
# empty cell under the new column, matching where Excel would land after
# entering the value in H6 and moving down.
$ws.Activate()
$ws.Range("H7").Select()
